# Add Login pattern with test-data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "no"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "password"

# --- Test rows ----------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "admin"
$ws.Range("A4").Value = 3
$ws.Range("C4").Value = "admin123"

# --- Alignment ----------------------------------------------------------
# "no" column (A) is centered with no border
$ws.Range("A1:A4").HorizontalAlignment = -4108   # xlCenter
# Header cells (B1:C1) are centered
$ws.Range("B1:C1").HorizontalAlignment = -4108   # xlCenter
# Data cells (B2:C4) are left aligned
$ws.Range("B2:C4").HorizontalAlignment = -4131   # xlLeft

# --- Borders --------------------------------------------------------------
# Thin box border around the username/password table (header + data)
$ws.Range("B1:C4").Borders.LineStyle = 1         # xlContinuous

# --- Column widths --------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 4.5
$ws.Columns.Item(2).ColumnWidth = 9.166666666666666
$ws.Columns.Item(3).ColumnWidth = 11.333333333333334

# --- Leave the selection where the original author left it ---------------
$ws.Range("D11").Select() | Out-Null
